$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 10

# Helper: write a literal text value into a cell. A leading apostrophe
# forces Excel to store the value verbatim as text (stops it from
# re-interpreting digit/date-shaped strings, and ensures a real - if
# empty - text cell gets created instead of no cell at all). Resetting
# the style back to "Normal" afterwards strips the quote-prefix flag
# that the apostrophe trick leaves behind, so the cell ends up with
# plain default formatting, matching a normal text cell.
function Set-TextCell($r, $c, $val) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Numeric columns
$ws.Cells.Item($row, 1).Value  = 131236926   # A  Id
$ws.Cells.Item($row, 2).Value  = 57881       # B  Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value  = 100049      # E  TaxonId
$ws.Cells.Item($row, 17).Value = 567473      # Q  Ost
$ws.Cells.Item($row, 18).Value = 6510086     # R  Nord
$ws.Cells.Item($row, 19).Value = 10          # S  Noggrannhet

# Plain text columns (no digit/date look-alikes, safe to assign directly)
$ws.Cells.Item($row, 4).Value  = "NT"                                    # D  Rödlistade
$ws.Cells.Item($row, 6).Value  = "Spillkråka"                            # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Dryocopus martius"                     # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(Linnaeus, 1758)"                      # H  Auktor
$ws.Cells.Item($row, 13).Value = "födosökande"                           # M  Aktivitet
$ws.Cells.Item($row, 16).Value = "Sydost Björkliden, Ög"                 # P  Lokalnamn
$ws.Cells.Item($row, 20).Value = "Östergötland"                          # T  Län
$ws.Cells.Item($row, 21).Value = "Norrköping"                            # U  Kommun
$ws.Cells.Item($row, 22).Value = "Östergötland"                          # V  Provins
$ws.Cells.Item($row, 23).Value = "Simonstorp"                            # W  Socken
$ws.Cells.Item($row, 29).Value = "Flera grova aspar lämpliga för bobygge i anslutning. Varav åtminstone två med äldre bohål."  # AC Publik kommentar
$ws.Cells.Item($row, 49).Value = "Anette Källman"                        # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Anette Källman"                        # AX Observatörer

# Text columns needing the forced-text helper: empty strings (must still
# create a present, empty text cell) and date-shaped strings (must not
# be auto-converted to a date serial number).
Set-TextCell $row 9  ""             # I  Antal (empty)
Set-TextCell $row 11 ""             # K  Ålder-Stadium (empty)
Set-TextCell $row 12 ""             # L  Kön (empty)
Set-TextCell $row 14 ""             # N  Metod (empty)
Set-TextCell $row 25 "2026-02-20"   # Y  Startdatum
Set-TextCell $row 27 "2026-02-20"   # AA Slutdatum
Set-TextCell $row 46 ""             # AT Bestämningsår (empty)
Set-TextCell $row 51 ""             # AY Projektnamn (empty)

# Boolean columns
$ws.Cells.Item($row, 30).Value = $false      # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false      # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false      # AG Ospontan
